# Update Efna3-Epha1 LR-pair sheet with new TPM-derived values
# Also reorders the MuSCs label (now appears before Efna3/Epha1 in the
# workbook shared strings) and appends 5 new rows for sending cluster = MuSCs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Epha1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.2420556666666667
$ws.Range("H2").Value2 = 0.726167
$ws.Range("I2").Value2 = 0.5314769098578004
$ws.Range("J2").Value2 = 0.5314769098578004
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.655411666666667
$ws.Range("N2").Value2 = 7.966235
$ws.Range("O2").Value2 = 0.1255433399118981
$ws.Range("P2").Value2 = 0.1255433399118982
$ws.Range("Q2").Value2 = 0.6427574412494445
$ws.Range("R2").Value2 = 5.784816971245
$ws.Range("S2").Value2 = 0.06672338634960308
$ws.Range("T2").Value2 = 0.06672338634960309

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna3"
$ws.Range("C3").Value = "Epha1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.2420556666666667
$ws.Range("H3").Value2 = 0.726167
$ws.Range("I3").Value2 = 0.5314769098578004
$ws.Range("J3").Value2 = 0.5314769098578004
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 7.831039666666666
$ws.Range("N3").Value2 = 23.493119
$ws.Range("O3").Value2 = 0.3702382146908386
$ws.Range("P3").Value2 = 0.3702382146908386
$ws.Range("Q3").Value2 = 1.895547527208111
$ws.Range("R3").Value2 = 17.059927744873
$ws.Range("S3").Value2 = 0.1967730622551558
$ws.Range("T3").Value2 = 0.1967730622551558

# Row 4: ECs -> Inflammatory-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Epha1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.2420556666666667
$ws.Range("H4").Value2 = 0.726167
$ws.Range("I4").Value2 = 0.5314769098578004
$ws.Range("J4").Value2 = 0.5314769098578004
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 6.127532
$ws.Range("N4").Value2 = 18.382596
$ws.Range("O4").Value2 = 0.2896992742608144
$ws.Range("P4").Value2 = 0.2896992742608145
$ws.Range("Q4").Value2 = 1.483203843281333
$ws.Range("R4").Value2 = 13.348834589532
$ws.Range("S4").Value2 = 0.1539684750721851
$ws.Range("T4").Value2 = 0.1539684750721851

# Row 5: ECs -> MuSCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna3"
$ws.Range("C5").Value = "Epha1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.2420556666666667
$ws.Range("H5").Value2 = 0.726167
$ws.Range("I5").Value2 = 0.5314769098578004
$ws.Range("J5").Value2 = 0.5314769098578004
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 2.108791333333333
$ws.Range("N5").Value2 = 6.326373999999999
$ws.Range("O5").Value2 = 0.09970006175963861
$ws.Range("P5").Value2 = 0.09970006175963862
$ws.Range("Q5").Value2 = 0.5104448920508888
$ws.Range("R5").Value2 = 4.594004028457999
$ws.Range("S5").Value2 = 0.05298828073664458
$ws.Range("T5").Value2 = 0.05298828073664458

# Row 6: ECs -> Resolving-Mac
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Efna3"
$ws.Range("C6").Value = "Epha1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 0.6666666666666666
$ws.Range("G6").Value2 = 0.2420556666666667
$ws.Range("H6").Value2 = 0.726167
$ws.Range("I6").Value2 = 0.5314769098578004
$ws.Range("J6").Value2 = 0.5314769098578004
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 2.428579666666666
$ws.Range("N6").Value2 = 7.285739
$ws.Range("O6").Value2 = 0.1148191093768101
$ws.Range("P6").Value2 = 0.1148191093768101
$ws.Range("Q6").Value2 = 0.587851470268111
$ws.Range("R6").Value2 = 5.290663232412999
$ws.Range("S6").Value2 = 0.06102370544421182
$ws.Range("T6").Value2 = 0.06102370544421183

# Row 7: FAPs -> ECs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna3"
$ws.Range("C7").Value = "Epha1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.2054156666666667
$ws.Range("H7").Value2 = 0.616247
$ws.Range("I7").Value2 = 0.4510271759376837
$ws.Range("J7").Value2 = 0.4510271759376837
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 2.655411666666667
$ws.Range("N7").Value2 = 7.966235
$ws.Range("O7").Value2 = 0.1255433399118981
$ws.Range("P7").Value2 = 0.1255433399118982
$ws.Range("Q7").Value2 = 0.5454631577827778
$ws.Range("R7").Value2 = 4.909168420045
$ws.Range("S7").Value2 = 0.05662345805824811
$ws.Range("T7").Value2 = 0.05662345805824812

# Row 8: FAPs -> FAPs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna3"
$ws.Range("C8").Value = "Epha1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value2 = 1
$ws.Range("F8").Value2 = 0.3333333333333333
$ws.Range("G8").Value2 = 0.2054156666666667
$ws.Range("H8").Value2 = 0.616247
$ws.Range("I8").Value2 = 0.4510271759376837
$ws.Range("J8").Value2 = 0.4510271759376837
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 7.831039666666666
$ws.Range("N8").Value2 = 23.493119
$ws.Range("O8").Value2 = 0.3702382146908386
$ws.Range("P8").Value2 = 0.3702382146908386
$ws.Range("Q8").Value2 = 1.608618233821444
$ws.Range("R8").Value2 = 14.477564104393
$ws.Range("S8").Value2 = 0.1669874963962188
$ws.Range("T8").Value2 = 0.1669874963962188

# Row 9: FAPs -> Inflammatory-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna3"
$ws.Range("C9").Value = "Epha1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value2 = 1
$ws.Range("F9").Value2 = 0.3333333333333333
$ws.Range("G9").Value2 = 0.2054156666666667
$ws.Range("H9").Value2 = 0.616247
$ws.Range("I9").Value2 = 0.4510271759376837
$ws.Range("J9").Value2 = 0.4510271759376837
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 6.127532
$ws.Range("N9").Value2 = 18.382596
$ws.Range("O9").Value2 = 0.2896992742608144
$ws.Range("P9").Value2 = 0.2896992742608145
$ws.Range("Q9").Value2 = 1.258691070801333
$ws.Range("R9").Value2 = 11.328219637212
$ws.Range("S9").Value2 = 0.1306622455410516
$ws.Range("T9").Value2 = 0.1306622455410517

# Row 10: FAPs -> MuSCs
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Efna3"
$ws.Range("C10").Value = "Epha1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value2 = 1
$ws.Range("F10").Value2 = 0.3333333333333333
$ws.Range("G10").Value2 = 0.2054156666666667
$ws.Range("H10").Value2 = 0.616247
$ws.Range("I10").Value2 = 0.4510271759376837
$ws.Range("J10").Value2 = 0.4510271759376837
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 2.108791333333333
$ws.Range("N10").Value2 = 6.326373999999999
$ws.Range("O10").Value2 = 0.09970006175963861
$ws.Range("P10").Value2 = 0.09970006175963862
$ws.Range("Q10").Value2 = 0.4331787775975555
$ws.Range("R10").Value2 = 3.898608998378
$ws.Range("S10").Value2 = 0.04496743729626245
$ws.Range("T10").Value2 = 0.04496743729626246

# Row 11: FAPs -> Resolving-Mac
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Efna3"
$ws.Range("C11").Value = "Epha1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value2 = 1
$ws.Range("F11").Value2 = 0.3333333333333333
$ws.Range("G11").Value2 = 0.2054156666666667
$ws.Range("H11").Value2 = 0.616247
$ws.Range("I11").Value2 = 0.4510271759376837
$ws.Range("J11").Value2 = 0.4510271759376837
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 2.428579666666666
$ws.Range("N11").Value2 = 7.285739
$ws.Range("O11").Value2 = 0.1148191093768101
$ws.Range("P11").Value2 = 0.1148191093768101
$ws.Range("Q11").Value2 = 0.4988683112814444
$ws.Range("R11").Value2 = 4.489814801533
$ws.Range("S11").Value2 = 0.05178653864590267
$ws.Range("T11").Value2 = 0.05178653864590268

# Row 12: MuSCs -> ECs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Efna3"
$ws.Range("C12").Value = "Epha1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value2 = 1
$ws.Range("F12").Value2 = 0.3333333333333333
$ws.Range("G12").Value2 = 0.007968333333333332
$ws.Range("H12").Value2 = 0.023905
$ws.Range("I12").Value2 = 0.01749591420451593
$ws.Range("J12").Value2 = 0.01749591420451593
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 2.655411666666667
$ws.Range("N12").Value2 = 7.966235
$ws.Range("O12").Value2 = 0.1255433399118981
$ws.Range("P12").Value2 = 0.1255433399118982
$ws.Range("Q12").Value2 = 0.02115920529722222
$ws.Range("R12").Value2 = 0.190432847675
$ws.Range("S12").Value2 = 0.00219649550404695
$ws.Range("T12").Value2 = 0.002196495504046951

# Row 13: MuSCs -> FAPs
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Efna3"
$ws.Range("C13").Value = "Epha1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value2 = 1
$ws.Range("F13").Value2 = 0.3333333333333333
$ws.Range("G13").Value2 = 0.007968333333333332
$ws.Range("H13").Value2 = 0.023905
$ws.Range("I13").Value2 = 0.01749591420451593
$ws.Range("J13").Value2 = 0.01749591420451593
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 7.831039666666666
$ws.Range("N13").Value2 = 23.493119
$ws.Range("O13").Value2 = 0.3702382146908386
$ws.Range("P13").Value2 = 0.3702382146908386
$ws.Range("Q13").Value2 = 0.06240033441055555
$ws.Range("R13").Value2 = 0.561603009695
$ws.Range("S13").Value2 = 0.006477656039464061
$ws.Range("T13").Value2 = 0.006477656039464061

# Row 14: MuSCs -> Inflammatory-Mac
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Efna3"
$ws.Range("C14").Value = "Epha1"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value2 = 1
$ws.Range("F14").Value2 = 0.3333333333333333
$ws.Range("G14").Value2 = 0.007968333333333332
$ws.Range("H14").Value2 = 0.023905
$ws.Range("I14").Value2 = 0.01749591420451593
$ws.Range("J14").Value2 = 0.01749591420451593
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 6.127532
$ws.Range("N14").Value2 = 18.382596
$ws.Range("O14").Value2 = 0.2896992742608144
$ws.Range("P14").Value2 = 0.2896992742608145
$ws.Range("Q14").Value2 = 0.04882621748666666
$ws.Range("R14").Value2 = 0.43943595738
$ws.Range("S14").Value2 = 0.005068553647577738
$ws.Range("T14").Value2 = 0.00506855364757774

# Row 15: MuSCs -> MuSCs
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Efna3"
$ws.Range("C15").Value = "Epha1"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value2 = 1
$ws.Range("F15").Value2 = 0.3333333333333333
$ws.Range("G15").Value2 = 0.007968333333333332
$ws.Range("H15").Value2 = 0.023905
$ws.Range("I15").Value2 = 0.01749591420451593
$ws.Range("J15").Value2 = 0.01749591420451593
$ws.Range("K15").Value2 = 3
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 2.108791333333333
$ws.Range("N15").Value2 = 6.326373999999999
$ws.Range("O15").Value2 = 0.09970006175963861
$ws.Range("P15").Value2 = 0.09970006175963862
$ws.Range("Q15").Value2 = 0.01680355227444444
$ws.Range("R15").Value2 = 0.15123197047
$ws.Range("S15").Value2 = 0.001744343726731576
$ws.Range("T15").Value2 = 0.001744343726731577

# Row 16: MuSCs -> Resolving-Mac
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Efna3"
$ws.Range("C16").Value = "Epha1"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value2 = 1
$ws.Range("F16").Value2 = 0.3333333333333333
$ws.Range("G16").Value2 = 0.007968333333333332
$ws.Range("H16").Value2 = 0.023905
$ws.Range("I16").Value2 = 0.01749591420451593
$ws.Range("J16").Value2 = 0.01749591420451593
$ws.Range("K16").Value2 = 3
$ws.Range("L16").Value2 = 1
$ws.Range("M16").Value2 = 2.428579666666666
$ws.Range("N16").Value2 = 7.285739
$ws.Range("O16").Value2 = 0.1148191093768101
$ws.Range("P16").Value2 = 0.1148191093768101
$ws.Range("Q16").Value2 = 0.01935173231055555
$ws.Range("R16").Value2 = 0.174165590795
$ws.Range("S16").Value2 = 0.002008865286695599
$ws.Range("T16").Value2 = 0.0020088652866956

Write-Output "Updated rows 2-16 with new TPM values"
